$wb = $excel.ActiveWorkbook

# Duplicate the "Login" sheet, place the copy right after "ProductList",
# and rename it to "AdminLogin" (this becomes the new Admin login page).
$login = $wb.Worksheets.Item("Login")
$productList = $wb.Worksheets.Item("ProductList")
$login.Copy($null, $productList)
$adminLogin = $wb.Worksheets.Item("Login (2)")
$adminLogin.Name = "AdminLogin"

# Login sheet: selection becomes the whole header range A1:B1 instead of D8.
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Range("A1:B1").Select()

# AdminLogin sheet: selection sits on F12, and this tab becomes the active one.
$adminLogin.Range("F12").Select()
$adminLogin.Activate()
